$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2778919.8
$ws.Range("J17").Value = 2778919.8
$ws.Range("L17").Value = 8336759.399999999
$ws.Range("N17").Value = -8337095.399999999
$ws.Range("H43").Value = 1456.5
$ws.Range("I43").Value = 1193.6666
$ws.Range("J43").Value = 1528.1818
$ws.Range("K43").Value = 1193.6666
$ws.Range("L43").Value = 1528.1818
$ws.Range("M43").Value = -1124.6666
$ws.Range("N43").Value = -1666.1818
$ws.Range("H45").Value = 2040
$ws.Range("I45").Value = 850
$ws.Range("K45").Value = 2550
$ws.Range("M45").Value = -2358
$ws.Range("H54").Value = 23999.166
$ws.Range("I54").Value = 9997.5
$ws.Range("J54").Value = 31000
$ws.Range("K54").Value = 9997.5
$ws.Range("L54").Value = 31000
$ws.Range("M54").Value = -9511.5
$ws.Range("N54").Value = -31972
$ws.Range("H63").Value = 33825
$ws.Range("J63").Value = 33825
$ws.Range("L63").Value = 33825
$ws.Range("N63").Value = -35073
$ws.Range("H66").Value = 33825
$ws.Range("J66").Value = 33825
$ws.Range("L66").Value = 101475
$ws.Range("N66").Value = -107715
$ws.Range("H75").Value = 111111
$ws.Range("J75").Value = 111111
$ws.Range("L75").Value = 111111
$ws.Range("N75").Value = -112983
$ws.Range("H78").Value = 111111
$ws.Range("J78").Value = 111111
$ws.Range("L78").Value = 333333
$ws.Range("N78").Value = -342693
$ws.Range("H137").Value = 2193.75
$ws.Range("I137").Value = 1564.6666
$ws.Range("J137").Value = 3137.375
$ws.Range("K137").Value = 4693.9998
$ws.Range("L137").Value = 9412.125
$ws.Range("M137").Value = -2143.9998
$ws.Range("N137").Value = -14512.125
$ws.Range("H138").Value = 2717.043
$ws.Range("I138").Value = 1646.5
$ws.Range("J138").Value = 3393.1755
$ws.Range("K138").Value = 4939.5
$ws.Range("L138").Value = 10179.5265
$ws.Range("M138").Value = 200.5
$ws.Range("N138").Value = -20459.5265

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1851.619
$ws.Range("I45").Value = 1834.3529
$ws.Range("J45").Value = 1925
$ws.Range("K45").Value = 1834.3529
$ws.Range("L45").Value = 1925
$ws.Range("M45").Value = -1457.3529
$ws.Range("N45").Value = -2679
$ws.Range("H74").Value = 10788.429
$ws.Range("I74").Value = 17574.857
$ws.Range("J74").Value = 4002
$ws.Range("K74").Value = 17574.857
$ws.Range("L74").Value = 4002
$ws.Range("M74").Value = -16700.857
$ws.Range("N74").Value = -5750
$ws.Range("H77").Value = 10788.429
$ws.Range("I77").Value = 17574.857
$ws.Range("J77").Value = 4002
$ws.Range("K77").Value = 87874.285
$ws.Range("L77").Value = 20010
$ws.Range("M77").Value = -83506.285
$ws.Range("N77").Value = -28746

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1585.3914
$ws.Range("I134").Value = 1521.1364
$ws.Range("K134").Value = 4563.4092
$ws.Range("M134").Value = -2028.4092

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3236.098
$ws.Range("I31").Value = 1176.2333
$ws.Range("J31").Value = 7098.3438
$ws.Range("K31").Value = 1176.2333
$ws.Range("L31").Value = 7098.3438
$ws.Range("M31").Value = -881.2333000000001
$ws.Range("N31").Value = -7688.3438
$ws.Range("H34").Value = 3236.098
$ws.Range("I34").Value = 1176.2333
$ws.Range("J34").Value = 7098.3438
$ws.Range("K34").Value = 1176.2333
$ws.Range("L34").Value = 7098.3438
$ws.Range("M34").Value = -974.2333000000001
$ws.Range("N34").Value = -7502.3438
$ws.Range("H58").Value = 2981.95
$ws.Range("I58").Value = 2465.818
$ws.Range("J58").Value = 3612.7778
$ws.Range("K58").Value = 2465.818
$ws.Range("L58").Value = 3612.7778
$ws.Range("M58").Value = -2262.818
$ws.Range("N58").Value = -4018.7778
$ws.Range("H134").Value = 2619.2693
$ws.Range("I134").Value = 2751.8572
$ws.Range("K134").Value = 8255.571599999999
$ws.Range("M134").Value = -5720.571599999999
$ws.Range("H136").Value = 2981.95
$ws.Range("I136").Value = 2465.818
$ws.Range("J136").Value = 3612.7778
$ws.Range("K136").Value = 7397.454000000001
$ws.Range("L136").Value = 10838.3334
$ws.Range("M136").Value = -4847.454000000001
$ws.Range("N136").Value = -15938.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1089.5555
$ws.Range("I68").Value = 1000.75
$ws.Range("K68").Value = 3002.25
$ws.Range("M68").Value = -2191.25
$ws.Range("H71").Value = 1089.5555
$ws.Range("I71").Value = 1000.75
$ws.Range("K71").Value = 9006.75
$ws.Range("M71").Value = -4950.75
$ws.Range("H127").Value = 980.5
$ws.Range("J127").Value = 980.5
$ws.Range("L127").Value = 2941.5
$ws.Range("N127").Value = -12861.5
$ws.Range("H131").Value = 857.9400000000001
$ws.Range("J131").Value = 857.9400000000001
$ws.Range("L131").Value = 2573.82
$ws.Range("N131").Value = -12653.82

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2588.923
$ws.Range("I132").Value = 2138.889
$ws.Range("J132").Value = 3601.5
$ws.Range("K132").Value = 6416.667
$ws.Range("L132").Value = 10804.5
$ws.Range("M132").Value = -3886.667
$ws.Range("N132").Value = -15864.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 25991.38
$ws.Range("I7").Value = 34070.613
$ws.Range("J7").Value = 3222.6365
$ws.Range("K7").Value = 34070.613
$ws.Range("L7").Value = 3222.6365
$ws.Range("M7").Value = -33958.613
$ws.Range("N7").Value = -3446.6365
$ws.Range("H22").Value = 1008
$ws.Range("J22").Value = 1340
$ws.Range("L22").Value = 1340
$ws.Range("N22").Value = -1930
$ws.Range("H27").Value = 1008
$ws.Range("J27").Value = 1340
$ws.Range("L27").Value = 1340
$ws.Range("N27").Value = -1554
$ws.Range("H126").Value = 25991.38
$ws.Range("I126").Value = 34070.613
$ws.Range("J126").Value = 3222.6365
$ws.Range("K126").Value = 102211.839
$ws.Range("L126").Value = 9667.9095
$ws.Range("M126").Value = -99741.83899999999
$ws.Range("N126").Value = -14607.9095
$ws.Range("H132").Value = 19542.908
$ws.Range("I132").Value = 19280.54
$ws.Range("K132").Value = 57841.62
$ws.Range("M132").Value = -55311.62
$ws.Range("H136").Value = 18577442
$ws.Range("I136").Value = 86675.75
$ws.Range("K136").Value = 260027.25
$ws.Range("M136").Value = -257477.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 361.93332
$ws.Range("I107").Value = 311
$ws.Range("J107").Value = 463.8
$ws.Range("K107").Value = 933
$ws.Range("L107").Value = 1391.4
$ws.Range("M107").Value = 987
$ws.Range("N107").Value = -5231.4
$ws.Range("H132").Value = 346513.3
$ws.Range("I132").Value = 477488.94
$ws.Range("J132").Value = 2702.25
$ws.Range("K132").Value = 1432466.82
$ws.Range("L132").Value = 8106.75
$ws.Range("M132").Value = -1429936.82
$ws.Range("N132").Value = -13166.75
